# GO4060_Coughlin.docx edit
#   1. Remove the "Justify" paragraph alignment from the abstract paragraph
#      (the one beginning "We seek to create a catalog of well-vetted
#      exoplanets ...").
#   2. Move the (hidden) _GoBack bookmark back to just after the word "is"
#      (its natural position after the preceding edit) and merge the run
#      split that used to straddle "false " / "positives" back into a
#      single run, by doing a Find/Replace over that phrase.

$d = $word.ActiveDocument

# --- 1. Un-justify the second paragraph -----------------------------------
$d.Paragraphs.Item(2).Alignment = 0   # wdAlignParagraphLeft

# --- 2. Relocate the _GoBack bookmark --------------------------------------
# Find the word "is" in "... over campaigns 4 and 5, which is expected ..."
# (search with surrounding spaces so we land on this occurrence, not the
# "is" substring inside another word).
$isRange = $d.Content
$isRange.Find.ClearFormatting()
$isRange.Find.Execute(" is ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterIsPos = $isRange.End - 1   # just past the "is", before the following space

# Drop the bookmark where it currently sits and re-add it right after "is"
if ($d.Bookmarks.Item("_GoBack").Start -ge 0) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackRange = $d.Range($afterIsPos, $afterIsPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 3. Re-merge the "false positives" runs --------------------------------
# A Find/Replace over text spanning the old run boundary collapses the two
# runs (". ... false " + "positives, ...") back into one contiguous run.
$mergeRange = $d.Content
$mergeRange.Find.ClearFormatting()
$mergeRange.Find.MatchWholeWord = $false
$mergeRange.Find.Execute("false positives", $false, $false, $false, $false, $false, $true, 1, $false, "false positives", 2) | Out-Null
